# This script reproduces the fix applied to the teacher-schedule workbook:
# for every day-of-week column (Lunes/Martes/Mierc/Jueves/Viernes, columns D-H)
# on each professor sheet, only the column(s) that actually correspond to the
# days the professor teaches keep their "x" mark; the spurious "x" marks that
# had been stamped across the whole row are cleared. The break-period row
# (row 9, "Receso") is cleared completely on every sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Miguel": keep x in columns E (Martes) and F (Mierc); clear D, G, H
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Miguel")
$ws.Activate() | Out-Null
$ws.Range("D5:D22").ClearContents() | Out-Null
$ws.Range("G5:G22").ClearContents() | Out-Null
$ws.Range("H5:H22").ClearContents() | Out-Null
$ws.Range("D9:H9").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Profesor2": keep x in columns D (Lunes) and E (Martes); clear F, G, H
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Profesor2")
$ws.Activate() | Out-Null
$ws.Range("F5:F22").ClearContents() | Out-Null
$ws.Range("G5:G22").ClearContents() | Out-Null
$ws.Range("H5:H22").ClearContents() | Out-Null
$ws.Range("D9:H9").ClearContents() | Out-Null
$ws.Range("D9").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Mauricio": keep x in column H (Viernes); clear D, E, F, G
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mauricio")
$ws.Activate() | Out-Null
$ws.Range("D5:D22").ClearContents() | Out-Null
$ws.Range("E5:E22").ClearContents() | Out-Null
$ws.Range("F5:F22").ClearContents() | Out-Null
$ws.Range("G5:G22").ClearContents() | Out-Null
$ws.Range("D9:H9").ClearContents() | Out-Null
$ws.Range("H9").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "profesor4": keep x in column H (Viernes); clear D, E, F, G
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("profesor4")
$ws.Activate() | Out-Null
$ws.Range("D5:D22").ClearContents() | Out-Null
$ws.Range("E5:E22").ClearContents() | Out-Null
$ws.Range("F5:F22").ClearContents() | Out-Null
$ws.Range("G5:G22").ClearContents() | Out-Null
$ws.Range("D9:H9").ClearContents() | Out-Null
$ws.Range("H5:H22").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Profesor5": keep x in column G (Jueves); clear D, E, F, H
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Profesor5")
$ws.Activate() | Out-Null
$ws.Range("D5:D22").ClearContents() | Out-Null
$ws.Range("E5:E22").ClearContents() | Out-Null
$ws.Range("F5:F22").ClearContents() | Out-Null
$ws.Range("H5:H22").ClearContents() | Out-Null
$ws.Range("D9:H9").ClearContents() | Out-Null
$ws.Range("G9").Select() | Out-Null

# ---------------------------------------------------------------------------
# Finish back on "Miguel" (the sheet that ends up active/selected), landing
# on F9 as the final selection.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Miguel")
$ws.Activate() | Out-Null
$ws.Range("F9").Select() | Out-Null
